# Applies the commit's data changes to the "Metadata" sheet:
#   - Experimental (row 7, column B) goes from blank to the literal text "false"
#   - Date (row 8, column B) is bumped to the new generation timestamp
#
# Everything else (Concepts sheet, styles, other metadata rows) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Experimental -> "false" -----------------------------------------------
# A plain assignment of the string "false" gets auto-coerced by Excel into the
# Boolean FALSE (same as typing `false` directly into a cell). We need the
# literal text "false" instead, so we build it in a scratch cell, trim it with
# a formula (which yields a genuine text result), and paste that text value
# (not the formula) into the target cell. This preserves the target cell's
# existing style/formatting.
$scratchA = $ws.Range("Z1")
$scratchB = $ws.Range("Z2")

$scratchA.Value = "falsex"
$scratchB.Formula = "=LEFT(Z1,5)"

$scratchB.Copy()
$ws.Range("B7").PasteSpecial(-4163)  # xlPasteValues

$scratchA.ClearContents()
$scratchB.ClearContents()

# --- Date --------------------------------------------------------------------
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"
